$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.931.20"
$ws.Range("E2").Value = "  -1.51%  "

# Row 3
$ws.Range("D3").Value = "1.868.57"
$ws.Range("E3").Value = "  -2.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.88"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5009"
$ws.Range("E7").Value = "  -1.88%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  -3.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08901"
$ws.Range("E9").Value = "  -8.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  -1.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.60"
$ws.Range("E11").Value = "  -1.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.398"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("E13").Value = "  -1.31%  "

# Row 14
$ws.Range("D14").Value = "1.868.27"
$ws.Range("E14").Value = "  -2.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.227"
$ws.Range("E15").Value = "  -2.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -2.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.11"
$ws.Range("E18").Value = "  -2.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06668"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.02"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.110"
$ws.Range("E22").Value = "  -2.04%  "

# Row 23
$ws.Range("D23").Value = "27.959.34"
$ws.Range("E23").Value = "  -1.74%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +0.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.271"
$ws.Range("E25").Value = "  -2.01%  "

# Row 26
$ws.Range("D26").Value = "2.081.97"
$ws.Range("E26").Value = "  -2.92%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.495"
$ws.Range("E27").Value = "  -6.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.26"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.64"
$ws.Range("E29").Value = "  -2.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.97"
$ws.Range("E30").Value = "  -1.30%  "

# Row 32
$ws.Range("E32").Value = "  -4.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.601"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  -0.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.507"
$ws.Range("E35").Value = "  -2.69%  "

# Row 36
$ws.Range("E36").Value = "  -2.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02393"
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2179"
$ws.Range("E38").Value = "  -1.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.283"
$ws.Range("E39").Value = "  +4.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").Value = "  -3.51%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6354"
$ws.Range("E41").Value = "  -0.82%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.46"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.903"
$ws.Range("E43").Value = "  -2.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.20"
$ws.Range("E45").Value = "  -3.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5997"
$ws.Range("E46").Value = "  -0.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.277"
$ws.Range("E47").Value = "  -0.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.664"
$ws.Range("E48").Value = "  -2.91%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.240"
$ws.Range("E49").Value = "  +3.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.991"
$ws.Range("E50").Value = "  -3.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.76"
$ws.Range("E51").Value = "  -2.66%  "
